$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.236.50'
$ws.Range("E2").Value = '  -2.77%  '
$ws.Range("D3").Value = '1.649.10'
$ws.Range("E3").Value = '  -3.33%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.44'
$ws.Range("E5").Value = '  -2.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3894'
$ws.Range("E7").Value = '  -1.54%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3879'
$ws.Range("E8").Value = '  -3.92%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.003'
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.369'
$ws.Range("E10").Value = '  -7.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '49.13'
$ws.Range("E11").Value = '  -6.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08487'
$ws.Range("E12").Value = '  -3.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.52'
$ws.Range("E13").Value = '  -5.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.164'
$ws.Range("E14").Value = '  -4.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001292'
$ws.Range("E15").Value = '  -4.86%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.525'
$ws.Range("E16").Value = '  -5.76%  '
$ws.Range("D17").Value = '1.649.73'
$ws.Range("E17").Value = '  -3.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.82'
$ws.Range("E18").Value = '  -1.60%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.20'
$ws.Range("E19").Value = '  +2.72%  '
$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06897'
$ws.Range("E20").Value = '  -3.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.978'
$ws.Range("E21").Value = '  -5.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.86'
$ws.Range("E23").Value = '  -4.21%  '
$ws.Range("D24").Value = '24.233.38'
$ws.Range("E24").Value = '  -2.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.369'
$ws.Range("E25").Value = '  +0.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.764'
$ws.Range("E26").Value = '  -7.34%  '
$ws.Range("E27").Value = '  -4.67%  '
$ws.Range("E28").Value = '  -2.02%  '
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.519'
$ws.Range("E29").Value = '  +0.58%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '143.01'
$ws.Range("E30").Value = '  -5.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.383'
$ws.Range("E31").Value = '  -12.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.464'
$ws.Range("E32").Value = '  -3.52%  '
$ws.Range("D33").Value = '1.828.83'
$ws.Range("E33").Value = '  -3.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08165'
$ws.Range("E34").Value = '  -4.48%  '
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.910'
$ws.Range("E35").Value = '  -3.81%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9963'
$ws.Range("E36").Value = '  -4.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02951'
$ws.Range("E37").Value = '  -6.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2732'
$ws.Range("E38").Value = '  -4.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09329'
$ws.Range("E39").Value = '  -2.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.481'
$ws.Range("E40").Value = '  +0.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.05'
$ws.Range("E41").Value = '  -7.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7674'
$ws.Range("E42").Value = '  -6.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.19'
$ws.Range("E43").Value = '  -5.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.12'
$ws.Range("E44").Value = '  -7.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.512'
$ws.Range("E45").Value = '  -6.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6920'
$ws.Range("E46").Value = '  -6.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.105'
$ws.Range("E47").Value = '  -3.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08481'
$ws.Range("E49").Value = '  -2.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.274'
$ws.Range("E50").Value = '  -7.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '134.49'
$ws.Range("E51").Value = '  -3.29%  '
